$d = $word.ActiveDocument

# --- Date / timestamp / runner-id updates in the title block and the
#     "generated on" sentence -------------------------------------------
$d.Content.Find.Execute("2021-07-03", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2021-08-25", 2)
$d.Content.Find.Execute("generated on 2021-08-25, 13:34:26 with", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "generated on 2021-08-25, 09:20:42 with", 2)
$d.Content.Find.Execute("runner-0277ea0f-project-18732201-concurrent-0", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "runner-fa6cab46-project-18732201-concurrent-0", 2)

# --- "Just a string" / "Just another string" become lead-ins for the
#     hyperlink that follows them, e.g. "Just a string: " -----------------
$d.Content.Find.Execute("Just a string", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Just a string: ", 2)
$d.Content.Find.Execute("Just another string", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Just another string: ", 2)

# Locate paragraphs by their text (Range.Text includes the trailing
# paragraph-mark "`r", so match against that too) instead of hard-coded
# indices, so the script stays correct even if paragraph numbering shifts.
# Re-scan fresh before each merge since deleting a paragraph mark shifts
# every later paragraph's index down by one.

# --- Merge "Just a string: " paragraph into the hyperlink paragraph that
#     follows it (delete the paragraph mark between them) -----------------
$count = $d.Paragraphs.Count
$idxLeadIn1 = -1
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq "Just a string: `r") { $idxLeadIn1 = $i }
}
$p1 = $d.Paragraphs($idxLeadIn1)
$mark1 = $p1.Range.End - 1
$d.Range($mark1, $p1.Range.End).Delete()

# --- Merge "Just another string: " paragraph into its hyperlink paragraph
$count = $d.Paragraphs.Count
$idxLeadIn2 = -1
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq "Just another string: `r") { $idxLeadIn2 = $i }
}
$p2 = $d.Paragraphs($idxLeadIn2)
$mark2 = $p2.Range.End - 1
$d.Range($mark2, $p2.Range.End).Delete()

# Re-locate the (now merged) paragraphs that hold each hyperlink so we can
# append the trailing " (Chapter N)" text right after the hyperlink run.
# Re-scan fresh before each insertion since InsertAfter shifts indices too.

# --- Append " (Chapter 1)" right after "Hello hyperlink!" ----------------
$count = $d.Paragraphs.Count
$idxMerged1 = -1
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq "Just a string: Hello hyperlink!`r") { $idxMerged1 = $i }
}
$p1b = $d.Paragraphs($idxMerged1)
$insert1 = $p1b.Range.End - 1
$d.Range($insert1, $insert1).InsertAfter(" (Chapter 1)")
# Touch the new run's formatting so Word materializes an (empty) rPr on it,
# matching the target markup (<w:rPr/> on the appended run).
$newRun1 = $d.Range($insert1, $insert1 + 12)
$newRun1.Font.Bold = $true
$newRun1.Font.Bold = $false

# --- Append " (Chapter 2)" right after "ARG's documentation" -------------
$count = $d.Paragraphs.Count
$idxMerged2 = -1
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq "Just another string: ARG's documentation`r") { $idxMerged2 = $i }
}
$p2b = $d.Paragraphs($idxMerged2)
$insert2 = $p2b.Range.End - 1
$d.Range($insert2, $insert2).InsertAfter(" (Chapter 2)")
$newRun2 = $d.Range($insert2, $insert2 + 12)
$newRun2.Font.Bold = $true
$newRun2.Font.Bold = $false
